$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rows 53-58 ("Chapter 10" detail rows): the trailing J column cell in each
#    of these rows is blank/unused and gets removed; the C/E/G numeric cells
#    keep their "00" number format but drop the (no-op, border-less) extra
#    formatting, and the D/F/H/I text cells drop back to the default style.
# ---------------------------------------------------------------------------
for ($r = 53; $r -le 58; $r++) {
    $ws.Range("C$r").NumberFormat = "00"
    $ws.Range("E$r").NumberFormat = "00"
    $ws.Range("G$r").NumberFormat = "00"
    $ws.Range("D$r").ClearFormats()
    $ws.Range("F$r").ClearFormats()
    $ws.Range("H$r").ClearFormats()
    $ws.Range("J$r").ClearFormats()
    $ws.Range("J$r").ClearContents()
}
# I53 loses its formatting entirely, while I54:I58 keep the "00" numeric xf
# (even though they hold text) just like the corresponding C/E/G cells.
$ws.Range("I53").ClearFormats()
for ($r = 54; $r -le 58; $r++) {
    $ws.Range("I$r").NumberFormat = "00"
}

# ---------------------------------------------------------------------------
# 2. New row 61: "Handling parent-child hierarchies" entry under chapter 11,
#    following directly below the existing row 60 ("Handling hierarchies").
#    Copy row 59's formatting (the last-row-of-section pattern, with the
#    bottom border) down onto row 61 so the same border/number-format xf
#    records get reused instead of new ones being minted.
# ---------------------------------------------------------------------------
$ws.Range("B59:K59").Copy() | Out-Null
$ws.Range("B61").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
# I61 needs the numeric "00" + bottom-border xf (same as C59/E59/G59) rather
# than the plain bottom-border-only xf the rest of row 59 used for column I.
$ws.Range("C59").Copy() | Out-Null
$ws.Range("I61").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

$ws.Range("C61").Value = 11
$ws.Range("D61").Value = "Handling hierarchies"
$ws.Range("E61").Value = 2
$ws.Range("F61").Value = "Handling parent-child hierarchies"
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = "Handling parent-child hierarchies"
$ws.Range("I61").Value = "PATH, PATHITEM, PATHLENGTH, LOOKUPVALUE, ISINSCOPE"
$ws.Range("B61").Formula = '=CONCAT(TEXT(C61,"00"),TEXT(E61,"00"),TEXT(G61,"00"))'

# ---------------------------------------------------------------------------
# 3. Refresh the selection so the active cell matches where the author left
#    off editing (column I of the newly-completed row).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("I53").Select() | Out-Null
